$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.219803940679071
$ws.Cells.Item(2, 4).Value = 0.04416460176246062
$ws.Cells.Item(2, 5).Value = 0.4396191691992186
$ws.Cells.Item(2, 6).Value = 1.528290297222981
$ws.Cells.Item(2, 7).Value = 1.673767081463382
$ws.Cells.Item(2, 8).Value = 1.010542361786406
$ws.Cells.Item(2, 9).Value = 3.961193417866525

$ws.Cells.Item(3, 2).Value = 1.069997499266549
$ws.Cells.Item(3, 4).Value = 0.03830832642090343
$ws.Cells.Item(3, 5).Value = 0.3823479279035951
$ws.Cells.Item(3, 6).Value = 1.363286054659909
$ws.Cells.Item(3, 7).Value = 1.472460773835678
$ws.Cells.Item(3, 8).Value = 0.9231199012414208
$ws.Cells.Item(3, 9).Value = 3.481483530427994

$ws.Cells.Item(4, 2).Value = 0.9775436500530645
$ws.Cells.Item(4, 4).Value = 0.03472637085036467
$ws.Cells.Item(4, 5).Value = 0.3474053453013823
$ws.Cells.Item(4, 6).Value = 1.262945067803372
$ws.Cells.Item(4, 7).Value = 1.349747030403535
$ws.Cells.Item(4, 8).Value = 0.8702008080870769
$ws.Cells.Item(4, 9).Value = 3.186307868167489

$ws.Cells.Item(5, 2).Value = 0.9397517899388959
$ws.Cells.Item(5, 4).Value = 0.03326964694228707
$ws.Cells.Item(5, 5).Value = 0.333216280922386
$ws.Cells.Item(5, 6).Value = 1.222288067622259
$ws.Cells.Item(5, 7).Value = 1.299949458491369
$ws.Cells.Item(5, 8).Value = 0.8488200623750686
$ws.Cells.Item(5, 9).Value = 3.06587693359225

$ws.Cells.Item(6, 2).Value = 0.9334695415141141
$ws.Cells.Item(6, 4).Value = 0.03302792226286044
$ws.Cells.Item(6, 5).Value = 0.3308630769826664
$ws.Cells.Item(6, 6).Value = 1.215550724464379
$ws.Cells.Item(6, 7).Value = 1.291692836908965
$ws.Cells.Item(6, 8).Value = 0.8452807431618226
$ws.Cells.Item(6, 9).Value = 3.04587117873993

$ws.Cells.Item(7, 2).Value = 0.9770344427033137
$ws.Cells.Item(7, 4).Value = 0.03470671369345268
$ws.Cells.Item(7, 5).Value = 0.3472137901143384
$ws.Cells.Item(7, 6).Value = 1.262395826126067
$ws.Cells.Item(7, 7).Value = 1.349074614279345
$ws.Cells.Item(7, 8).Value = 0.8699117225502562
$ws.Cells.Item(7, 9).Value = 3.184684259837752

$ws.Cells.Item(8, 2).Value = 1.168250105382356
$ws.Cells.Item(8, 4).Value = 0.04214216942487781
$ws.Cells.Item(8, 5).Value = 0.4198225367166089
$ws.Cells.Item(8, 6).Value = 1.47118785694326
$ws.Cells.Item(8, 7).Value = 1.604163365212003
$ws.Cells.Item(8, 8).Value = 0.98023762357343
$ws.Cells.Item(8, 9).Value = 3.795928923952232

$ws.Cells.Item(9, 2).Value = 1.539389405565316
$ws.Cells.Item(9, 4).Value = 0.05685744970638495
$ws.Cells.Item(9, 5).Value = 0.5642327170640442
$ws.Cells.Item(9, 6).Value = 1.888902018709359
$ws.Cells.Item(9, 7).Value = 2.112127782647462
$ws.Cells.Item(9, 8).Value = 1.202918706961327
$ws.Cells.Item(9, 9).Value = 4.989014804984947

$ws.Cells.Item(10, 2).Value = 1.809644099351601
$ws.Cells.Item(10, 4).Value = 0.06778656051733378
$ws.Cells.Item(10, 5).Value = 0.6719446234360191
$ws.Cells.Item(10, 6).Value = 2.201659492699889
$ws.Cells.Item(10, 7).Value = 2.491042104895939
$ws.Cells.Item(10, 8).Value = 1.370841485233257
$ws.Cells.Item(10, 9).Value = 5.861568061238756

$ws.Cells.Item(11, 2).Value = 1.932050532541837
$ws.Cells.Item(11, 4).Value = 0.07279217337988086
$ws.Cells.Item(11, 5).Value = 0.7213816670721798
$ws.Cells.Item(11, 6).Value = 2.345399064368451
$ws.Cells.Item(11, 7).Value = 2.664885521213535
$ws.Cells.Item(11, 8).Value = 1.44827547428082
$ws.Cells.Item(11, 9).Value = 6.25753114434076

$ws.Cells.Item(12, 2).Value = 1.978324344475823
$ws.Cells.Item(12, 4).Value = 0.07469321267035411
$ws.Cells.Item(12, 5).Value = 0.7401723167001393
$ws.Cells.Item(12, 6).Value = 2.400055170825851
$ws.Cells.Item(12, 7).Value = 2.730945692543116
$ws.Cells.Item(12, 8).Value = 1.477756451588164
$ws.Cells.Item(12, 9).Value = 6.407322560668433

$ws.Cells.Item(13, 2).Value = 1.968361989861023
$ws.Cells.Item(13, 4).Value = 0.0742835337146488
$ws.Cells.Item(13, 5).Value = 0.7361221926845474
$ws.Cells.Item(13, 6).Value = 2.388273756458403
$ws.Cells.Item(13, 7).Value = 2.71670796254466
$ws.Cells.Item(13, 8).Value = 1.471400018306724
$ws.Cells.Item(13, 9).Value = 6.375069214559801

$ws.Cells.Item(14, 2).Value = 1.935859098657602
$ws.Cells.Item(14, 4).Value = 0.0729484585700817
$ws.Cells.Item(14, 5).Value = 0.7229261437163217
$ws.Cells.Item(14, 6).Value = 2.349891060477688
$ws.Cells.Item(14, 7).Value = 2.670315630292464
$ws.Cells.Item(14, 8).Value = 1.450697669861711
$ws.Cells.Item(14, 9).Value = 6.269857673436093

$ws.Cells.Item(15, 2).Value = 1.915939811074622
$ws.Cells.Item(15, 4).Value = 0.07213142522000737
$ws.Cells.Item(15, 5).Value = 0.7148524956228215
$ws.Cells.Item(15, 6).Value = 2.326410295331868
$ws.Cells.Item(15, 7).Value = 2.641929398215893
$ws.Cells.Item(15, 8).Value = 1.438037758433723
$ws.Cells.Item(15, 9).Value = 6.205392518134602

$ws.Cells.Item(16, 2).Value = 1.801633628445018
$ws.Cells.Item(16, 4).Value = 0.06746016902452823
$ws.Cells.Item(16, 5).Value = 0.6687232001523142
$ws.Cells.Item(16, 6).Value = 2.192296469740057
$ws.Cells.Item(16, 7).Value = 2.479712166869206
$ws.Cells.Item(16, 8).Value = 1.365802710240757
$ws.Cells.Item(16, 9).Value = 5.835670441680008

$ws.Cells.Item(17, 2).Value = 1.731372246121737
$ws.Cells.Item(17, 4).Value = 0.06460364233217319
$ws.Cells.Item(17, 5).Value = 0.6405414693056741
$ws.Cells.Item(17, 6).Value = 2.110407240501445
$ws.Cells.Item(17, 7).Value = 2.38058685249996
$ws.Cells.Item(17, 8).Value = 1.321762217021217
$ws.Cells.Item(17, 9).Value = 5.608601612490247

$ws.Cells.Item(18, 2).Value = 1.690909632146656
$ws.Cells.Item(18, 4).Value = 0.06296379515380579
$ws.Cells.Item(18, 5).Value = 0.6243729184997164
$ws.Cells.Item(18, 6).Value = 2.063443689084153
$ws.Cells.Item(18, 7).Value = 2.323710192705676
$ws.Cells.Item(18, 8).Value = 1.2965291561203
$ws.Cells.Item(18, 9).Value = 5.477907426105105

$ws.Cells.Item(19, 2).Value = 1.677201143000843
$ws.Cells.Item(19, 4).Value = 0.06240909297251562
$ws.Cells.Item(19, 5).Value = 0.6189053326469747
$ws.Cells.Item(19, 6).Value = 2.047565729167701
$ws.Cells.Item(19, 7).Value = 2.304475839262807
$ws.Cells.Item(19, 8).Value = 1.288002240429705
$ws.Cells.Item(19, 9).Value = 5.433641513037287

$ws.Cells.Item(20, 2).Value = 1.73885689351971
$ws.Cells.Item(20, 4).Value = 0.06490739369498044
$ws.Cells.Item(20, 5).Value = 0.6435371872360207
$ws.Cells.Item(20, 6).Value = 2.119110209985365
$ws.Cells.Item(20, 7).Value = 2.391124544869683
$ws.Cells.Item(20, 8).Value = 1.326440215915113
$ws.Cells.Item(20, 9).Value = 5.632782903627174

$ws.Cells.Item(21, 2).Value = 1.945408141630935
$ws.Cells.Item(21, 4).Value = 0.07334044737054057
$ws.Cells.Item(21, 5).Value = 0.726800189790751
$ws.Cells.Item(21, 6).Value = 2.361158763589458
$ws.Cells.Item(21, 7).Value = 2.683935813114317
$ws.Cells.Item(21, 8).Value = 1.45677408627796
$ws.Cells.Item(21, 9).Value = 6.30076505272973

$ws.Cells.Item(22, 2).Value = 2.079940855283951
$ws.Cells.Item(22, 4).Value = 0.07888443913425647
$ws.Cells.Item(22, 5).Value = 0.7816282444840539
$ws.Cells.Item(22, 6).Value = 2.520670277776958
$ws.Cells.Item(22, 7).Value = 2.876650817821144
$ws.Cells.Item(22, 8).Value = 1.542882076632281
$ws.Cells.Item(22, 9).Value = 6.736443306619606

$ws.Cells.Item(23, 2).Value = 2.00818100542449
$ws.Cells.Item(23, 4).Value = 0.07592231310724173
$ws.Cells.Item(23, 5).Value = 0.752325552932561
$ws.Cells.Item(23, 6).Value = 2.435410421870188
$ws.Cells.Item(23, 7).Value = 2.773666093207737
$ws.Cells.Item(23, 8).Value = 1.496837023020134
$ws.Cells.Item(23, 9).Value = 6.503998686294892

$ws.Cells.Item(24, 2).Value = 1.735473297350211
$ws.Cells.Item(24, 4).Value = 0.06477006025150445
$ws.Cells.Item(24, 5).Value = 0.6421827193271525
$ws.Cells.Item(24, 6).Value = 2.115175238321626
$ws.Cells.Item(24, 7).Value = 2.386360107228029
$ws.Cells.Item(24, 8).Value = 0.98023762357343
$ws.Cells.Item(24, 9).Value = 5.621851005971394

$ws.Cells.Item(25, 2).Value = 1.439406031741441
$ws.Cells.Item(25, 4).Value = 0.05285854641996934
$ws.Cells.Item(25, 5).Value = 0.524910646920091
$ws.Cells.Item(25, 6).Value = 1.774930254878313
$ws.Cells.Item(25, 7).Value = 1.973780395381823
$ws.Cells.Item(25, 8).Value = 1.14195305260256
$ws.Cells.Item(25, 9).Value = 4.666914883151094
